# Applies the "1_0_yearly_repeated" dictionary edit:
#  - child_id (row 3) valueType changes from "integer" to "text", and its
#    unit (previously "numeric") is cleared.
#  - a number of rows that previously had no "unit" value (meta-variables
#    such as age_years, whe_, asthma_, asthma_med_, URTI_, LRTI_, eczema_,
#    rash_, rash_loc_, ...) get "numeric" filled in as their unit.
#  - sheet "Variables" becomes the active/selected sheet (with C56
#    selected), and sheet "Categories" is left with its whole second row
#    selected and is no longer the active tab.

$wb = $excel.ActiveWorkbook

$wsVariables = $wb.Worksheets.Item("Variables")
$wsCategories = $wb.Worksheets.Item("Categories")

# --- child_id (row 3): valueType integer -> text, unit numeric -> (blank)
$wsVariables.Range("B3").Value = "text"
$wsVariables.Range("C3").Value = ""

# --- fill in the previously-blank "unit" column with "numeric" for the
#     rows that describe meta-variables / plain numeric variables.
$numericUnitRows = @(4, 5, 6, 7, 8, 9, 10, 19, 20, 22, 23, 24, 25, 26, 31, 32, 33)
foreach ($r in $numericUnitRows) {
    $wsVariables.Range("C$r").Value = "numeric"
}

# --- selection / active sheet bookkeeping, matching the saved workbook state
$wsCategories.Activate()
$wsCategories.Rows.Item(2).Select()

$wsVariables.Activate()
$wsVariables.Range("C56").Select()
